# Update the single data row (row 2) of the report so that it reflects the
# new fixture values. Plain textual replacements just use .Value, but the
# "Prix"/"Quantite" amounts and the "Date" column are numeric/date-looking
# strings that must stay as literal text (they were stored as shared text
# strings, not real numbers/dates). Pre-formatting those cells as Text ("@")
# before assigning the value prevents Excel from silently converting them
# into a number / date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name
$ws.Range("A2").Value = "shopping with Kahlo"

# Date (keep as literal text, not an actual date value)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1950-11-30"

# Magasin/Restaurant
$ws.Range("E2").Value = "Mediocre Wool Coat"

# Produit
$ws.Range("F2").Value = "Celery Seed"

# Prix (keep as literal text, not a real number)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5.0"

# Quantite (keep as literal text, not a real number)
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "5.0"
